$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$grey = 14277081   # RGB(217,217,217) - "White, darker 15%"
$red  = 255        # RGB(255,0,0)
$black = 0         # RGB(0,0,0)

# --- Insert a new row at the top for the "PETUNJUK" instructions banner ---
$ws.Rows("1:1").Insert()

# --- Register the bold/red font used by the rich-text run in the styles
#     table (mirrors what real Excel does when you format part of a cell's
#     text via the Font dialog). Use a scratch cell well away from the used
#     range, then clear it so it leaves no trace in the saved sheet. ---
$scratch = $ws.Range("Z100")
$scratch.Font.Bold = $true
$scratch.Font.Color = $red
$scratch.Clear()

# --- Row 1: grey banner row with instructions ---
$ws.Range("A1:B1").Interior.Color = $grey

$ws.Range("C1").Value = "PETUNJUK" + [char]10 + "SILAHKAN MASUKKAN KODE GURU" + [char]10 + "DI ANTARA SYMBOL [ ]"
$ws.Range("C1").Interior.Color = $grey
$ws.Range("C1").HorizontalAlignment = -4108
$ws.Range("C1").VerticalAlignment = -4108
$ws.Range("C1").WrapText = $true

$title = $ws.Range("C1").Characters(1, 8)
$title.Font.Bold = $true
$title.Font.Color = $red

$body = $ws.Range("C1").Characters(9, 52)
$body.Font.Name = "Calibri"
$body.Font.Size = 11
$body.Font.Color = $black

$ws.Rows("1").RowHeight = 87

# --- Row 2 keeps the original header row content (now shifted down) ---
$ws.Rows("2").RowHeight = 59.4

# --- Row 3 (previously row 2): code placeholder becomes a literal "[1]" ---
$ws.Range("C3").Value = "[1]"

# --- Widen column C to fit the new instructions ---
$ws.Columns("C").ColumnWidth = 29.77734375

# --- Misc view state ---
$ws.Range("F2").Select()
$ws.PageSetup.Orientation = 1

Write-Output "done"
